$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Common Word"
$ws.Range("B1").Value = "Total Frequency"
$ws.Range("C1").Value = "Websites"

# Data rows
$ws.Range("A2").Value = 'Cookies'
$ws.Range("B2").Value = 83
$ws.Range("C2").Value = 'https://www.100-dakar.com (14), https://www.benjaminwahl.at (18), https://www.drehorgelkabarett.at (14), https://www.ottosaxinger.at (3), https://www.peligro.at (14), https://www.schuledesungehorsams.at (2), https://www.skodone.at (18)'
$ws.Range("A3").Value = 'Page'
$ws.Range("B3").Value = 35
$ws.Range("C3").Value = 'https://www.hungaromedia.at (8), https://www.kuenstlerinnen.at (8), https://www.luckeneder-art.at (8), https://www.platform-socialism.org (3), https://www.regional-express.org (8)'
$ws.Range("A4").Value = 'Linz'
$ws.Range("B4").Value = 52
$ws.Range("C4").Value = 'https://www.freie-medien.at (10), https://www.freizeitundkommunikation.at (3), https://www.linzfmr.at (18), https://www.steingeschichten.at (21)'
$ws.Range("A5").Value = 'March'
$ws.Range("B5").Value = 49
$ws.Range("C5").Value = 'https://www.das-kollektiv.at (8), https://www.feminismus-krawall.at (16), https://www.fiftitu.at (19), https://www.unkraut-comics.at (6)'
$ws.Range("A6").Value = 'Art'
$ws.Range("B6").Value = 29
$ws.Range("C6").Value = 'https://www.eipcp.net (13), https://www.kairus.org (6), https://www.negentropy-sport.net (2), https://www.radical-openness.org (8)'
$ws.Range("A7").Value = 'School'
$ws.Range("B7").Value = 294
$ws.Range("C7").Value = 'https://www.alteschule-gutau.at (4), https://www.derschueler.at (5), https://www.die-schule.at (285)'
$ws.Range("A8").Value = 'Culture'
$ws.Range("B8").Value = 37
$ws.Range("C8").Value = 'https://www.frauenkultur.at (17), https://www.igkultur.at (15), https://www.interregnum.live (5)'
$ws.Range("A9").Value = 'Uhr'
$ws.Range("B9").Value = 33
$ws.Range("C9").Value = 'https://www.fro.at (21), https://www.rudolfhabringer.at (12)'
$ws.Range("A10").Value = 'Film'
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 'https://www.corpushomini.info (3), https://www.doublehappiness.at (12)'
$ws.Range("A11").Value = 'Magdalena'
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 'https://www.magdalenareiter.at (2), https://www.themagdalenaproject.org (10)'
$ws.Range("A12").Value = 'Casino'
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 'https://www.frf.at (6), https://www.photosalonhelga.com (4)'
$ws.Range("A13").Value = 'Andreas'
$ws.Range("B13").Value = 8
$ws.Range("C13").Value = 'https://www.andreaskurz.net (2), https://www.andreaszingerle.com (6)'

# Header styling: bold, centered, top-aligned, thin box border
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108  # xlCenter
$header.VerticalAlignment = -4160    # xlTop
$header.Borders.LineStyle = 1        # xlContinuous
$header.Borders.Weight = 2           # xlThin

